$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = '19 hours ago ... Iran''s Raisi reiterates warnings as Israel mulls response to air attack. As the world calls for calm, President Raisi vows ''slightest attack'' will be met with a ...'

# Row 3
$ws.Range("A3").Value = 'LIVE: Pakistan vs New Zealand – T20 international cricket | Cricket ...'
$ws.Range("B3").Value = 'Published On 18 Apr 2024'
$ws.Range("C3").Value = '5 hours ago ... The sides will face each other five times as they continue their warm-up to the T20 World Cup. ... Hello and welcome to our live coverage of ...'
$ws.Range("D3").Value = 'LIVE:Pakistanvs'
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("A4").Value = 'Nigeria | Today''s latest from Al Jazeera'
$ws.Range("C4").Value = '19 hours ago ... Why is Germany maintaining economic ties with China? German Chancellor Olaf Scholz has been on a three-day visit to China in a bid to shore up economic ties.'
$ws.Range("D4").Value = 'NigeriaToday''s'

# Row 5
$ws.Range("A5").Value = 'Opinion | Today''s latest from Al Jazeera'
$ws.Range("B5").Value = 'No Date'
$ws.Range("C5").Value = '19 hours ago ... The conflict in Sudan has displaced over two million people, triggering one of the world''s largest humanitarian crises. Opinion by Amitabh BeharAmitabh ...'
$ws.Range("D5").Value = 'OpinionToday''s'
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("A6").Value = 'Elections | Today''s latest from Al Jazeera'
$ws.Range("C6").Value = '19 hours ago ... Croatians vote in election pitting the PM against the country''s president ... Exit polls expected minutes after voting ends at 7pm (17:00 GMT), with official ...'
$ws.Range("D6").Value = 'ElectionsToday''s'

# Row 7
$ws.Range("A7").Value = 'Video | Today''s latest from Al Jazeera'
$ws.Range("B7").Value = 'No Date'
$ws.Range("C7").Value = '19 hours ago ... Dominican FM on Haiti gang violence crisis: Spillover threat? Roberto Alvarez Gil, Dominican Republic FM, discusses the effect of Haiti''s criminal gang turmoil.'
$ws.Range("D7").Value = 'VideoToday''s'
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("A8").Value = 'Listen Live to Al Jazeera | Al Jazeera'
$ws.Range("C8").Value = '2 days ago ... Live Broadcast. NEWS 30min. 2:00PM - 2:30PM. Up-to-date news and analysis from around the world.'
$ws.Range("D8").Value = 'ListenLiveto'

# Row 9
$ws.Range("A9").Value = 'Philippines | Today''s latest from Al Jazeera'
$ws.Range("B9").Value = 'No Date'
$ws.Range("C9").Value = '2 days ago ... Filipino migrant workers in European country allege wage theft, salary deductions and passport confiscation. Published On 17 Apr 202417 Apr 2024.'
$ws.Range("D9").Value = 'PhilippinesToday''s'
$ws.Range("E9").Value = 0

# Row 10
$ws.Range("A10").Value = 'Turkey | Today''s latest from Al Jazeera'
$ws.Range("B10").Value = 'No Date'
$ws.Range("C10").Value = '3 days ago ... Stay on top of Turkey latest developments on the ground with Al Jazeera''s fact-based news, exclusive video footage, photos and updated maps.'
$ws.Range("D10").Value = 'TurkeyToday''s'

# Row 11
$ws.Range("A11").Value = 'South Africa | Today''s latest from Al Jazeera'
$ws.Range("B11").Value = 'No Date'
$ws.Range("C11").Value = '3 days ago ... Former President Jacob Zuma hopes to run for office for the opposition uMkhonto weSizwe Party (MK) in May elections. Published On 12 Apr 202412 Apr 2024.'
$ws.Range("D11").Value = 'SouthAfrica'
